# Tidsplan_Vinkaelder.xlsx update
# "opdateret tidsplan med programmerings del" - update the schedule with the
# programming/UI-design section: rename the old "Programmering" header,
# add four new task columns (Brugerflade minimum/selvvalgte ure, Backend
# vaerdiproduktion, Newsfeed), re-color the progress cells and widen the
# new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: rename column L, add headers for the 4 new columns
# ---------------------------------------------------------------------
$ws.Range("L1").Value = "Program brugerfladedesign "
$ws.Range("M1").Value = "Brugerflade, minimum"
$ws.Range("N1").Value = "Brugerflade, selvvalgte ure"
$ws.Range("O1").Value = "Backend, v" + [char]0x00E6 + "rdiproduktion"
$ws.Range("P1").Value = "Newsfeed"

# ---------------------------------------------------------------------
# 2. Column widths (L..Q) for the new / widened columns
# ---------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 31.25
$ws.Columns.Item(13).ColumnWidth = 27.59
$ws.Columns.Item(14).ColumnWidth = 25.92
$ws.Columns.Item(15).ColumnWidth = 25.09
$ws.Columns.Item(16).ColumnWidth = 13.59
$ws.Columns.Item(17).ColumnWidth = 21.25

# ---------------------------------------------------------------------
# 3. Cell fills for the progress tracker
#    Green  = RGB(0,176,80)   (existing "done" fill)
#    Yellow = RGB(255,255,0)  (existing "in progress" fill)
#    Red    = RGB(192,0,0)    (existing "not started" fill)
# ---------------------------------------------------------------------
$green  = 5287936   # RGB(0,176,80)
$yellow = 65535      # RGB(255,255,0)
$red    = 192        # RGB(192,0,0)

# Row 3
$ws.Range("L3").ClearFormats()

# Row 4
$ws.Range("L4").Interior.Color = $green
$ws.Range("O4").ClearFormats()

# Row 5
$ws.Range("J5").Interior.Color = $yellow
$ws.Range("L5").Interior.Color = $yellow
$ws.Range("M5").Interior.Color = $yellow
$ws.Range("O5").Interior.Color = $yellow

# Row 6 - L6 loses its yellow fill, M6 gets green text on yellow fill
$ws.Range("L6").ClearFormats()
$ws.Range("M6").Interior.Color = $yellow
$ws.Range("M6").Font.Color = $green
$ws.Range("O6").Interior.Color = $yellow

# Row 7
$ws.Range("H7").Interior.Color = $yellow
$ws.Range("I7").Interior.Color = $yellow
$ws.Range("L7").ClearFormats()
$ws.Range("M7").Interior.Color = $green
$ws.Range("N7").Interior.Color = $yellow
$ws.Range("O7").Interior.Color = $green
$ws.Range("P7").Interior.Color = $yellow

# Row 8
$ws.Range("L8").ClearFormats()
$ws.Range("M8").Interior.Color = $green
$ws.Range("N8").Interior.Color = $yellow
$ws.Range("O8").Interior.Color = $green
$ws.Range("P8").Interior.Color = $yellow

# Row 9
$ws.Range("L9").ClearFormats()
$ws.Range("M9").Interior.Color = $red
$ws.Range("N9").Interior.Color = $yellow
$ws.Range("O9").Interior.Color = $red
$ws.Range("P9").Interior.Color = $yellow

# Row 10
$ws.Range("L10").ClearFormats()
$ws.Range("M10").Interior.Color = $red
$ws.Range("N10").Interior.Color = $yellow
$ws.Range("O10").Interior.Color = $red
$ws.Range("P10").Interior.Color = $yellow

# ---------------------------------------------------------------------
# 4. Row 11 - drop the now-unused leftover formatted cells
# ---------------------------------------------------------------------
$ws.Range("K11").Clear()
$ws.Range("L11").Clear()

# ---------------------------------------------------------------------
# 5. View: scroll the frozen pane over and move the selection
# ---------------------------------------------------------------------
$ws.Range("F1").Select()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("C18").Select()
